$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T for rows 2-16.
# (F and L remain unchanged at 1.)
$data = @{
    2  = @{ E=3; G=1.384145666666667; H=4.152437; I=0.1014617184198512; J=0.1334061399754118; K=3; M=21.084959; N=63.25487699999999; O=0.03381320693734752; P=0.03509122472428063; Q=29.18465463169433; R=262.661891685249; S=0.003430746081149312; T=0.004681384837476012 }
    3  = @{ E=3; G=1.384145666666667; H=4.152437; I=0.1014617184198512; J=0.1334061399754118; K=3; M=81.06331633333333; N=243.189949; O=0.12999838843446; P=0.1349118606466557; Q=112.2034380284126; R=1009.830942255713; S=0.01318985988237162; T=0.01799807056577099 }
    4  = @{ E=3; G=1.384145666666667; H=4.152437; I=0.1014617184198512; J=0.1334061399754118; K=3; M=136.9994176666667; N=410.998253; O=0.2197011461990087; P=0.2280050605000741; Q=189.6271502991735; R=1706.644352692561; S=0.02229125583216238; T=0.03041727501617511 }
    5  = @{ E=3; G=1.384145666666667; H=4.152437; I=0.1014617184198512; J=0.1334061399754118; K=3; M=316.292811; N=948.878433; O=0.5072276531881493; P=0.5263990368430604; Q=437.795323743469; R=3940.157913691221; S=0.05146418932253793; T=0.07022486359200726 }
    6  = @{ E=3; G=1.384145666666667; H=4.152437; I=0.1014617184198512; J=0.1334061399754118; K=2; M=68.131198; N=136.262396; O=0.1092596052410345; P=0.07559281728592908; Q=94.30350247650867; R=565.821014859052; S=0.01108566730162994; T=0.01008454596398238 }
    7  = @{ E=3; G=2.458038666666667; H=7.374116000000001; I=0.1801810554109116; J=0.2369096391566985; K=3; M=21.084959; N=63.25487699999999; O=0.03381320693734752; P=0.03509122472428063; Q=51.82764450708133; R=466.448800563732; S=0.006092499312798832; T=0.008313449386995941 }
    8  = @{ E=3; G=2.458038666666667; H=7.374116000000001; I=0.1801810554109116; J=0.2369096391566985; K=3; M=81.06331633333333; N=243.189949; O=0.12999838843446; P=0.1349118606466557; Q=199.2567659955649; R=1793.310893960084; S=0.02342324682983864; T=0.03196192022375799 }
    9  = @{ E=3; G=2.458038666666667; H=7.374116000000001; I=0.1801810554109116; J=0.2369096391566985; K=3; M=136.9994176666667; N=410.998253; O=0.2197011461990087; P=0.2280050605000741; Q=336.7498659354831; R=3030.748793419348; S=0.03958598439712437; T=0.05401659660897377 }
    10 = @{ E=3; G=2.458038666666667; H=7.374116000000001; I=0.1801810554109116; J=0.2369096391566985; K=3; M=316.292811; N=948.878433; O=0.5072276531881493; P=0.5263990368430604; Q=777.459959426692; R=6997.139634840229; S=0.09139281388504057; T=0.1247090058709231 }
    11 = @{ E=3; G=2.458038666666667; H=7.374116000000001; I=0.1801810554109116; J=0.2369096391566985; K=2; M=68.131198; N=136.262396; O=0.1092596052410345; P=0.07559281728592908; Q=167.4691190903227; R=1004.814714541936; S=0.01968651098610916; T=0.0179086670660477 }
    12 = @{ E=2; G=9.799863999999999; H=19.599728; I=0.7183572261692373; J=0.6296842208678898; K=3; M=21.084959; N=63.25487699999999; O=0.03381320693734752; P=0.03509122472428063; Q=206.629730645576; R=1239.778383873456; S=0.02428996154339937; T=0.02209639049980868 }
    13 = @{ E=2; G=9.799863999999999; H=19.599728; I=0.7183572261692373; J=0.6296842208678898; K=3; M=81.06331633333333; N=243.189949; O=0.12999838843446; P=0.1349118606466557; Q=794.4094754556453; R=4766.456852733872; S=0.09338528172224973; T=0.08495186985712669 }
    14 = @{ E=2; G=9.799863999999999; H=19.599728; I=0.7183572261692373; J=0.6296842208678898; K=3; M=136.9994176666667; N=410.998253; O=0.2197011461990087; P=0.2280050605000741; Q=1342.575661212531; R=8055.453967275183; S=0.1578239059697219; T=0.1435711888749252 }
    15 = @{ E=2; G=9.799863999999999; H=19.599728; I=0.7183572261692373; J=0.6296842208678898; K=3; M=316.292811; N=948.878433; O=0.5072276531881493; P=0.5263990368430604; Q=3099.626531977704; R=18597.75919186622; S=0.3643706499805708; T=0.3314651673801301 }
    16 = @{ E=2; G=9.799863999999999; H=19.599728; I=0.7183572261692373; J=0.6296842208678898; K=2; M=68.131198; N=136.262396; O=0.1092596052410345; P=0.07559281728592908; Q=667.676474557072; R=2670.705898228288; S=0.0784874269532954; T=0.047599604255899 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
